$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overal Stats")
$ws.Range("CJ1").Value = 43982
